$wb = $excel.ActiveWorkbook

# Grab the two existing worksheets by their current (pre-edit) tab order.
$wsTabelle1 = $wb.Worksheets.Item(1)   # currently "Tabelle1", holds "Region" in A1
$wsSheet1   = $wb.Worksheets.Item(2)   # currently "Sheet1", empty

# --- Build the "modelrun_setting" sheet (was "Sheet1") ---------------------
$wsSheet1.Name = "modelrun_setting"

$wsSheet1.Range("A1").Value = "module"
$wsSheet1.Range("B1").Value = "is_activated "
$wsSheet1.Range("C1").Value = "explanatory_comment"
$wsSheet1.Range("A1:D1").Font.Bold = $true

$wsSheet1.Range("A2").Value = "A"
$wsSheet1.Range("B2").Value = $true
$wsSheet1.Range("C2").Value = "vehicle stock"

$wsSheet1.Range("A3").Value = "B"
$wsSheet1.Range("B3").Value = $true
$wsSheet1.Range("C3").Value = "regional"

$wsSheet1.Range("A4").Value = "C"
$wsSheet1.Range("B4").Value = $false
$wsSheet1.Range("C4").Value = "interregional"

$wsSheet1.Range("A5").Value = "D"
$wsSheet1.Range("B5").Value = $false
$wsSheet1.Range("C5").Value = "mode infrastructure"

$wsSheet1.Range("A6").Value = "E"
$wsSheet1.Range("B6").Value = $false
$wsSheet1.Range("C6").Value = "fueling infrastructure"

$wsSheet1.Range("A7").Value = "F"
$wsSheet1.Range("B7").Value = $false
$wsSheet1.Range("C7").Value = "fuel supply infrastructure"

$wsSheet1.Columns.Item(2).ColumnWidth = 11.7

$wsSheet1.PageSetup.PaperSize = 9
$wsSheet1.PageSetup.Orientation = 1

[void]$wsSheet1.Range("H16").Select()

# --- Build the "set_nb_names" sheet (was "Tabelle1") ------------------------
$wsTabelle1.Range("A1").Value = "set_type"
$wsTabelle1.Range("B1").Value = "index"
$wsTabelle1.Range("C1").Value = "name"

$wsTabelle1.Name = "set_nb_names"
[void]$wsTabelle1.Range("A2").Select()

# --- Reorder the tabs: "modelrun_setting" goes in front of "set_nb_names" --
$wsSheet1.Move($wsTabelle1)
